$d = $word.ActiveDocument

# NOTE: we deliberately do NOT pass a replacement string into Find.Execute's
# own Replace parameter. Doing so causes this host to aggressively coalesce
# every same-formatted run in the paragraph into a single run, destroying the
# original run/rsid structure. Instead we locate each match with a
# search-only Find (Replace = wdReplaceNone = 0) and then assign the new text
# directly onto the matched Range, which only rewrites that one run's text
# and preserves every sibling run untouched - exactly matching the diff.

# 1) Body text: "A TERE," -> "A QWER,"  (bold run "TERE")
$bodyRng = $d.Content
$null = $bodyRng.Find.Execute("TERE", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$bodyRng.Text = "QWER"

# 2) Header: work through the primary header of the (only) section.
$hdr = $d.Sections.Item(1).Headers.Item(1)

# Helper range that we will keep sliding forward through the header story
$rng = $hdr.Range.Duplicate
$rng.Start = 0
$rng.End = $hdr.Range.End

# "DIRETORIA DE ENSINO REGIAO TRE" -> "...QWER"
$null = $rng.Find.Execute("TRE", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$rng.Text = "QWER"
$rng.Collapse(0)
$rng.End = $hdr.Range.End

# "TERE - DEP." -> "QWER - DEP."
$null = $rng.Find.Execute("TERE", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$rng.Text = "QWER"
$rng.Collapse(0)
$rng.End = $hdr.Range.End

# "Tre, nº Tre - Tre - Tre - Tre" -> "Qwer, nº Qwer - Qewr - Qewr - Qwer"
$treReplacements = @("Qwer", "Qwer", "Qewr", "Qewr", "Qwer")
for ($i = 0; $i -lt $treReplacements.Length; $i++) {
    $null = $rng.Find.Execute("Tre", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    $rng.Text = $treReplacements[$i]
    $rng.Collapse(0)
    $rng.End = $hdr.Range.End
}

# "CEP: tre    -    Tel: tre" -> "CEP: qwer    -    Tel: qwer"
for ($i = 0; $i -lt 2; $i++) {
    $null = $rng.Find.Execute("tre", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
    $rng.Text = "qwer"
    $rng.Collapse(0)
    $rng.End = $hdr.Range.End
}

# "Email: tre" -> "Email: qwer"
$null = $rng.Find.Execute("tre", $true, $false, $false, $false, $false, $true, 1, $false, $null, 0)
$rng.Text = "qwer"
$rng.Collapse(0)
$rng.End = $hdr.Range.End

Write-Host "Done."
